# Fruta / hortaliza, semanal
# Insert a new week's worth of data (3 rows) for "Terminal La Palmera de La
# Serena - Platano" right above the existing block that starts at row 448,
# shifting all the following rows down by 3 (dimension grows from T484 to
# T487).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 448 (existing rows 448:484 move to 451:487)
$ws.Rows("448:450").Insert()

# Row 448: Pinton
$ws.Range("A448").Value = 8
$ws.Range("B448").Value = "Terminal La Palmera de La Serena"
$ws.Range("C448").Value = "Coquimbo"
$ws.Range("D448").Value = 44578
$ws.Range("D448").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E448").Value = 4
$ws.Range("F448").Value = "Fruta"
$ws.Range("G448").Value = 100108
$ws.Range("H448").Value = "Tropicales y subtropicales"
$ws.Range("I448").Value = 100108006
$ws.Range("J448").Value = "Plátano"
$ws.Range("K448").Value = "Sin especificar"
$ws.Range("L448").Value = "Pintón"
$ws.Range("M448").Value = 80
$ws.Range("N448").Value = 15000
$ws.Range("O448").Value = 15000
$ws.Range("P448").Value = 15000
$ws.Range("Q448").Value = "$/caja 20 kilos"
$ws.Range("R448").Value = "Ecuador"
$ws.Range("S448").Value = 750
$ws.Range("T448").Value = 20

# Row 449: Primera Maduro
$ws.Range("A449").Value = 8
$ws.Range("B449").Value = "Terminal La Palmera de La Serena"
$ws.Range("C449").Value = "Coquimbo"
$ws.Range("D449").Value = 44578
$ws.Range("D449").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E449").Value = 4
$ws.Range("F449").Value = "Fruta"
$ws.Range("G449").Value = 100108
$ws.Range("H449").Value = "Tropicales y subtropicales"
$ws.Range("I449").Value = 100108006
$ws.Range("J449").Value = "Plátano"
$ws.Range("K449").Value = "Sin especificar"
$ws.Range("L449").Value = "Primera Maduro"
$ws.Range("M449").Value = 120
$ws.Range("N449").Value = 17000
$ws.Range("O449").Value = 17000
$ws.Range("P449").Value = 17000
$ws.Range("Q449").Value = "$/caja 20 kilos"
$ws.Range("R449").Value = "Ecuador"
$ws.Range("S449").Value = 850
$ws.Range("T449").Value = 20

# Row 450: Primera Pinton
$ws.Range("A450").Value = 8
$ws.Range("B450").Value = "Terminal La Palmera de La Serena"
$ws.Range("C450").Value = "Coquimbo"
$ws.Range("D450").Value = 44578
$ws.Range("D450").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E450").Value = 4
$ws.Range("F450").Value = "Fruta"
$ws.Range("G450").Value = 100108
$ws.Range("H450").Value = "Tropicales y subtropicales"
$ws.Range("I450").Value = 100108006
$ws.Range("J450").Value = "Plátano"
$ws.Range("K450").Value = "Sin especificar"
$ws.Range("L450").Value = "Primera Pintón"
$ws.Range("M450").Value = 120
$ws.Range("N450").Value = 18000
$ws.Range("O450").Value = 18000
$ws.Range("P450").Value = 18000
$ws.Range("Q450").Value = "$/caja 20 kilos"
$ws.Range("R450").Value = "Ecuador"
$ws.Range("S450").Value = 900
$ws.Range("T450").Value = 20
